# Updated cryptos list — applies the row-shift + value refresh described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (one decimal point) need an explicit Text
# number format first, otherwise Excel COM auto-converts the assigned string into
# a floating-point number and silently drops significant trailing/leading zeros.
$textCells = @(
    "D4", "D5", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16",
    "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D25", "D27", "D28",
    "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38",
    "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48",
    "D49", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the refreshed values (coin name, link, price, 1h volume-change columns).
$ws.Range("D2").Value = "29.495.65"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.901.03"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.50%  "
$ws.Range("D5").Value = "338.68"
$ws.Range("E5").Value = "  +4.07%  "
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  -1.49%  "
$ws.Range("D8").Value = "0.4010"
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "47.15"
$ws.Range("E9").Value = "  -2.07%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.08048"
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "0.9915"
$ws.Range("E11").Value = "  -2.18%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "23.24"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.918.97"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.945"
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "7.100"
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "89.08"
$ws.Range("E16").Value = "  -2.86%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "0.06814"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("B18").Value = "BinanceUSD"
$ws.Range("C18").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D18").Value = "1.005"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.00001019"
$ws.Range("E19").Value = "  -1.98%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "17.34"
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "29.521.70"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "5.516"
$ws.Range("E23").Value = "  -2.72%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "11.63"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "2.153"
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "2.132.43"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "157.27"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "6.502"
$ws.Range("E28").Value = "  -2.21%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "19.61"
$ws.Range("E29").Value = "  -2.01%  "
$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Value = "2.055"
$ws.Range("E30").Value = "  -2.95%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "119.30"
$ws.Range("E31").Value = "  -1.30%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "0.9955"
$ws.Range("E32").Value = "  -2.03%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").Value = "0.09539"
$ws.Range("E33").Value = "  -0.97%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "5.484"
$ws.Range("E34").Value = "  -2.84%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.387"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "3.530"
$ws.Range("E36").Value = "  -0.48%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "0.06409"
$ws.Range("E37").Value = "  +4.92%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.02247"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "1.199"
$ws.Range("E39").Value = "  +1.42%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.5828"
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "10.55"
$ws.Range("E41").Value = "  -3.21%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "7.742"
$ws.Range("E42").Value = "  -4.03%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.1820"
$ws.Range("E43").Value = "  -1.48%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "2.443"
$ws.Range("E44").Value = "  +2.06%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "1.266"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "12.15"
$ws.Range("E46").Value = "  -2.45%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.07362"
$ws.Range("E47").Value = "  -3.16%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "0.5491"
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "1.956"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "116.39"
$ws.Range("E50").Value = "  -2.09%  "
$ws.Range("D51").Value = "71.20"
$ws.Range("E51").Value = "  -1.52%  "
